$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 417
    3  = 418
    4  = 421
    5  = 424
    6  = 426
    7  = 427
    8  = 429
    9  = 430
    10 = 432
    11 = 434
    12 = 437
    13 = 438
    14 = 441
    15 = 12
    16 = 44
    17 = 98
    18 = 137
    19 = 159
    20 = 208
    21 = 240
    22 = 295
    23 = 353
    24 = 388
    25 = 411
    26 = 472
    27 = 517
}

foreach ($row in $values.Keys) {
    $ws.Range("A$row").Value = $values[$row]
}
